$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5516.6335,
    5387.312,
    5339.824000000001,
    5306.819000000001,
    5338.83,
    5422.714499999999,
    5875.107,
    7214.4275,
    8467.018,
    13640.6165,
    15763.02,
    15072.0395,
    15259.6675,
    15815.744,
    15955.3275,
    16114.0735,
    16384.389,
    15932.637,
    15036.9065,
    13330.5025,
    11770.8115,
    9312.075499999999,
    6773.7985,
    5542.130999999999
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
